# Update res_line/loading_percent values for the 380 kV case (rows 2-25, A2:A25 = 0..23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @(2, 4, 5, 6, 7, 10, 11, 12)  # B, D, E, F, G, J, K, L

$newValues = @{
    2 = @(19.40724839967093, 4.224014819138913, 10.49168018273442, 67.00831813771507, 3.816963889247042, 10.8198374715759, 17.40621676252509, 11.57093747002975)
    3 = @(19.42037544281899, 4.093170547565819, 10.48530197573399, 66.16903701268269, 3.821033700799329, 10.80287517787516, 17.41054528819653, 11.62335430146027)
    4 = @(19.4336601485137, 4.010937432710511, 10.48136693515264, 65.65229266315593, 3.823659902145234, 10.79227442009194, 17.41982035980609, 11.65858724392999)
    5 = @(19.44038689404548, 3.97699063797426, 10.47975830135992, 65.44150869475952, 3.824762244961114, 10.78790740823387, 17.42526098685786, 11.67371013712301)
    6 = @(19.44158314164091, 3.971328691817858, 10.47949087250419, 65.40650011389411, 3.824947233366021, 10.7871794121149, 17.42626460054326, 11.67626745426251)
    7 = @(19.43374555242482, 4.010481324838621, 10.48134526151479, 65.64945059373292, 3.823674638406028, 10.7922157163005, 17.41988701399906, 11.65878810035205)
    8 = @(19.41069002878286, 4.179316000028145, 10.48948466723992, 66.7193200099875, 3.818340812507023, 10.81402721018953, 17.40633441799278, 11.58837714066772)
    9 = @(19.40693546581859, 4.493765153543601, 10.50531317618705, 68.79942040413599, 3.808885554909869, 10.85534336488878, 17.43233320639701, 11.47455765923322)
    10 = @(19.42941278664656, 4.712723450725488, 10.51688492408191, 70.30796458615491, 3.802542833986239, 10.88483770706924, 17.48348547198573, 11.40581112408624)
    11 = @(19.4450949548717, 4.809367408061481, 10.52214191503644, 70.98813341409425, 3.79978677017717, 10.89807436943839, 17.51368351060297, 11.37778585200262)
    12 = @(19.45181486302899, 4.845513201718778, 10.5241319934465, 71.24467201514307, 3.798761573855, 10.90306120998918, 17.52610982173931, 11.3676420094427)
    13 = @(19.4503329179951, 4.837748998889063, 10.5237034163059, 71.18946988767283, 3.798981548948985, 10.90198833747059, 17.52338962145805, 11.36980579448939)
    14 = @(19.4456321668217, 4.812350321890483, 10.5223056503876, 71.00926043533066, 3.799702057296989, 10.89848515078452, 17.51468600640712, 11.37694190859398)
    15 = @(19.44285447237228, 4.796733437298304, 10.5214494078655, 70.89873877975332, 3.800145790509695, 10.89633602483632, 17.50948365880646, 11.38137407512576)
    16 = @(19.42849746529295, 4.706345722168819, 10.51654125815503, 70.26337854039878, 3.802725539917529, 10.88396905906857, 17.4816509582066, 11.4077081411207)
    17 = @(19.42108548950992, 4.650118410884396, 10.51352873238121, 69.87193954926163, 3.804341153958276, 10.87633635892258, 17.46634776596648, 11.4246962977594)
    18 = @(19.41733621442895, 4.617500552851348, 10.51179531986722, 69.6462312700337, 3.805282587481749, 10.87192922489147, 17.45819834915772, 11.43477304219024)
    19 = @(19.41615512513146, 4.606409881679074, 10.51120828649468, 69.56971836300012, 3.805603435552744, 10.87043411325406, 17.45555131458767, 11.4382372912421)
    20 = @(19.42182134217039, 4.656132805980966, 10.51384948926085, 69.91366827853172, 3.804167909912143, 10.87715063043998, 17.46790931179294, 11.42285623818758)
    21 = @(19.44699171431607, 4.819822962138629, 10.52271622213018, 71.06222138010592, 3.799489926374285, 10.89951481351054, 17.5172156248894, 11.37483312470713)
    22 = @(19.46799459515807, 4.924164295783247, 10.52850765472316, 71.80681551004992, 3.796540165912105, 10.91398206676786, 17.55521270296482, 11.34618006333189)
    23 = @(19.45636969269763, 4.868724560532751, 10.52541684913896, 71.41001483012374, 3.798104706571387, 10.90627414221382, 17.53440692791388, 11.36122213948038)
    24 = @(19.42148706827257, 4.653414605606089, 10.5137044796873, 69.89480477905738, 3.804246194303952, 10.87678255747292, 17.46720131631716, 11.42368716358284)
    25 = @(19.40351386907567, 4.410690208395749, 10.5010429508532, 68.23952850721616, 3.811336780912182, 10.8443193203421, 17.41966528446622, 11.50274276419983)
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Cells.Item($row, $columns[$i]).Value = $values[$i]
    }
}
